# Trade #27 closed at 2026-02-16 22:54:47 - base_strategy DOWN +0.000%
#
# Appends a new row (row 28) describing trade #27 to both the "All Trades"
# and "base_strategy" worksheets - mirroring the row already present for
# trade #26 (row 27) on each sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 28

    # --- Trade # --------------------------------------------------------
    $ws.Cells.Item($row, 1).Value = 27

    # --- Date / Time ------------------------------------------------------
    # These look like dates/times, so Excel would normally auto-convert them
    # to date/time serial numbers on assignment. Force the cells to Text
    # first so the literal strings are preserved, then clear the formatting
    # back to the default (General) so no stray number format lingers on
    # the cell - matching the plain, unstyled text cells used elsewhere in
    # this sheet.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-16"
    $ws.Cells.Item($row, 2).ClearFormats()

    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = "22:54:46"
    $ws.Cells.Item($row, 3).ClearFormats()

    # --- Strategy / Side / Status / Reason text fields --------------------
    $ws.Cells.Item($row, 4).Value = "base_strategy"
    $ws.Cells.Item($row, 5).Value = "DOWN"

    # --- Entry price --------------------------------------------------
    $ws.Cells.Item($row, 6).Value = 49.999998

    # --- Exit price (blank until the trade closes) -------------------
    $ws.Cells.Item($row, 7).Value = ""

    # --- Status ------------------------------------------------------
    $ws.Cells.Item($row, 8).Value = "OPEN"

    # --- P&L %, P&L $, Capital After, slippage, confidence -------------
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6

    # --- Entry reason --------------------------------------------------
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"

    # --- Exit reason (blank until the trade closes) -------------------
    $ws.Cells.Item($row, 16).Value = ""

    # --- Duration (min) -------------------------------------------------
    $ws.Cells.Item($row, 17).Value = 0
}
